$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5: new test case IPA0001 ---
# Copy formatting from existing bordered row (row 4) to the new row first
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A5:E5").PasteSpecial(-4122) | Out-Null

$ws.Range("A5").Value = 'IPA0001'
$ws.Range("B5").Value = "placeholder-b5"
$ws.Range("C5").Value = 'Verify that, accessing of the URL  takes the user to DRA application Landing page || Verify that DRA Landing page, displays application branding and logo || Verify that DRA Landing page, contains feature promotion and iconography in the marketing section || Verify that DRA Landing page, displays link to privacy statement and terms of use. || verify that DRA Landing page, displays the message and email id on the DRA landing page "Having trouble with sign-in? please contact DRA_support@thomsonreuters.com "||'
$ws.Range("D5").Value = "Y"
$ws.Range("E5").Value = ""

# Wrap text for the long-text columns (matches new cellXfs style w/ wrapText)
$ws.Range("B5").WrapText = $true
$ws.Range("C5").WrapText = $true

# Hyperlink on B5 -> Jira URL; display attribute = the URL itself
$ws.Hyperlinks.Add($ws.Range("B5"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4176", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-4176") | Out-Null
$ws.Range("B5").Value = 'OPQA-4176||OPQA-4178||OPQA-4179||OPQA-4182||OPQA-4187||OPQA-4189 '

$ws.Rows.Item(5).RowHeight = 120

# --- Row 6: new test case IPA0002 ---
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A6:E6").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = 'IPA0002'
$ws.Range("B6").Value = 'OPQA-4249 ||OPQA-4247 ||OPQA-4238'
$ws.Range("C6").Value = ' Verify that when linking a social with a matching email, if the user click [X] cross mark on the screen then he will be taken back to the DRA Login page. || Verify that text on the modal "Already have an account? .. ||Verify that when linking a social with a matching email, if the user clicks outside the Linking modal on the screen then nothing should happens'
$ws.Range("D6").Value = "Y"
$ws.Range("E6").Value = ""

$ws.Range("B6").WrapText = $true
$ws.Range("C6").WrapText = $true

$ws.Rows.Item(6).RowHeight = 90

$ws.Range("C6").Select() | Out-Null

$wb.Application.CutCopyMode = $false
